$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Dividende je Aktie"
$ws.Range("D1").Value = "Zahlungen pro Jahr"
$ws.Range("A1").Value = "Name"
$ws.Range("B18").Value = 60

$ws.Range("B10").Select()
